$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:G11 (naive component forecaster bug fix - rows shifted down
# by one period, with a new leading row of near-zero error stats, and the old
# trailing row dropped).
$data = @(
    @([double]"2.554709391686116E-07", [double]"8.01143367166098E-07", [double]"3.360900671017877E-12", [double]"1.833275939682261E-06", [double]"1.879105582334942E-06", 15),
    @([double]"-0.005035358036557557", [double]"0.3637990656943072", [double]"0.1855351867275108", [double]"0.4307379559865961", [double]"0.4469673578304695", 14),
    @([double]"-0.01898232632975465", [double]"0.3165751929851393", [double]"0.1363897707384557", [double]"0.3693098573534908", [double]"0.3838817902180699", 13),
    @([double]"0.00780077805212256", [double]"0.3794346132818944", [double]"0.1847968220473418", [double]"0.4298800088947401", [double]"0.4489210943938488", 12),
    @([double]"0.03532231998103826", [double]"0.2564852750112934", [double]"0.1004647918617584", [double]"0.3169618145167623", [double]"0.3303616777566779", 11),
    @([double]"0.002000156080238219", [double]"0.3449448151542904", [double]"0.1826303857083173", [double]"0.4273527649475516", [double]"0.4504644332784307", 10),
    @([double]"-0.05875859174690282", [double]"0.3248866822092542", [double]"0.1283271032654252", [double]"0.3582277254281489", [double]"0.3748117494910135", 9),
    @([double]"-0.02790046359007027", [double]"0.3120137538489823", [double]"0.1287774372279015", [double]"0.3588557331684997", [double]"0.3824716629792676", 8),
    @([double]"-0.03188629724616485", [double]"0.2657389849834738", [double]"0.1191175486163901", [double]"0.3451341023665874", [double]"0.3711930514458828", 7),
    @([double]"-0.0881379738596985", [double]"0.4006321887415465", [double]"0.2481252247171507", [double]"0.4981216966938408", [double]"0.5370552175868598", 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $col = 2 + $c
        $ws.Cells.Item($row, $col).Value = $vals[$c]
    }
}
